# Add payment 79174445 (Cash) 2025-08-18T08:51:52
#
# 1) Row 16, col A ("phone") was stored as a text value "79174445"; it
#    should become a genuine number 79174445 (matches every other "phone"
#    cell in the sheet).
# 2) A new row 17 is appended with the new payment record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix A16: text "79174445" -> number 79174445 ------------------
$ws.Range("A16").Value = 79174445

# --- 2) Append row 17 --------------------------------------------------
# Phone number is written as text (matches the diff's inlineStr <t>79174445</t>).
# A leading apostrophe forces Excel to store it as text instead of a number;
# resetting the style afterwards drops the quote-prefix formatting Excel
# applies so the cell ends up plain text with the sheet's default style.
$ws.Range("A17").Value = "'79174445"
$ws.Range("A17").Style = "Normal"

# B17/F17 ("amount"/"discount_applied") are blank in this record, stored
# as empty (quote-prefixed) text cells, matching the sheet's existing
# convention for not-applicable numeric fields.
$ws.Range("B17").Value = "'"
$ws.Range("B17").Style = "Normal"

$ws.Range("C17").Value = "Cash"
$ws.Range("D17").Value = "2025-08-18T08:51:52"
$ws.Range("E17").Value = 20

$ws.Range("F17").Value = "'"
$ws.Range("F17").Style = "Normal"

$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 20

Write-Output "done"
